$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.308.98'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.679.00'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.24'
$ws.Range('D5').Style = 'Normal'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5260'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.93%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.008'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.28%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06455'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.93'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.71%  '
$ws.Range('E10').Style = 'Normal'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07515'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.28%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.681.78'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.524'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.35%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5806'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.04%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008492'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.76'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.340.27'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.87'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.23'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.48%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.209'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.13'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.834'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.74%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1251'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.68%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.80'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.25%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06485'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.69%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.95%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.84%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('B31').Style = 'Normal'

$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C31').Style = 'Normal'

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.594'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.15%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('B32').Style = 'Normal'

$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C32').Style = 'Normal'

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.595'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.50%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.666'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.31%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.029'
$ws.Range('D34').Style = 'Normal'

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.38%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6236'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.54%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.91%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.37%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.426'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +4.08%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01623'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.09%  '
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.108.25'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.74%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8770'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.04%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.67%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.46'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.43%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.830.39'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.94%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('B45').Style = 'Normal'

$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('C45').Style = 'Normal'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000110'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.55%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'Aave'
$ws.Range('B46').Style = 'Normal'

$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('C46').Style = 'Normal'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.92'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.35%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.233'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.35%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.41%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.30%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4295'
$ws.Range('D50').Style = 'Normal'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.079'
$ws.Range('D51').Style = 'Normal'

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.37%  '
$ws.Range('E51').Style = 'Normal'

